# Updated cryptos list on Thu Oct 26 13:14:49 UTC 2023 with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) figures pulled from coinranking.com,
# and fix the row ordering for a couple of coins whose rank swapped (B/C/D/E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "34.209.63"
$ws.Range("E2").Value = "  -0.62%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.818.24"
$ws.Range("E3").Value = "  +1.75%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.34"
$ws.Range("E5").Value = "  +0.26%  "

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.558"
$ws.Range("E6").Value = "  +0.96%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.10%  "

# Row 8: Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.13"
$ws.Range("E8").Value = "  -4.26%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  +3.92%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0739"
$ws.Range("E10").Value = "  +11.94%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("E11").Value = "  -0.14%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.080.51"
$ws.Range("E12").Value = "  +1.81%  "

# Row 13: Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.06"
$ws.Range("E13").Value = "  -0.22%  "

# Row 14: WrappedEther
$ws.Range("D14").Value = "1.816.21"
$ws.Range("E14").Value = "  +1.74%  "

# Row 15: Polygon
$ws.Range("E15").Value = "  +1.53%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "34.215.37"
$ws.Range("E16").Value = "  -0.45%  "

# Row 17: Polkadot
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.34"
$ws.Range("E17").Value = "  +2.32%  "

# Row 18: Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.73"
$ws.Range("E18").Value = "  +0.81%  "

# Row 19: BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "250.13"
$ws.Range("E19").Value = "  -2.36%  "

# Row 20: ShibaInu
$ws.Range("E20").Value = "  +7.97%  "

# Row 21: Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.12"
$ws.Range("E21").Value = "  +6.46%  "

# Row 22: Dai
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.12%  "

# Row 23: Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.28"
$ws.Range("E23").Value = "  +1.90%  "

# Row 24: Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +0.65%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.70"
$ws.Range("E25").Value = "  +2.05%  "

# Row 26: EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.70"
$ws.Range("E26").Value = "  +1.64%  "

# Row 27: Cosmos
$ws.Range("E27").Value = "  +2.77%  "

# Row 28: Stellar
$ws.Range("E28").Value = "  +0.99%  "

# Row 29: BinanceUSD
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30: Hedera
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0534"
$ws.Range("E30").Value = "  +3.77%  "

# Row 31: Filecoin
$ws.Range("E31").Value = "  +0.43%  "

# Row 32: PancakeSwap
$ws.Range("E32").Value = "  +2.50%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.61"
$ws.Range("E33").Value = "  +0.84%  "

# Row 34: LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("E34").Value = "  -1.10%  "

# Row 35: Maker
$ws.Range("D35").Value = "1.433.69"
$ws.Range("E35").Value = "  -0.79%  "

# Row 36: TrustWalletToken
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.644"
$ws.Range("E36").Value = "  +2.80%  "

# Row 37: ImmutableX
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").Value = "  +0.63%  "

# Row 38: VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0191"
$ws.Range("E38").Value = "  +0.96%  "

# Row 39: ARBITRUM
$ws.Range("E39").Value = "  +8.05%  "

# Row 40: Aave
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.93"
$ws.Range("E40").Value = "  -1.47%  "

# Row 41: MXToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  -3.54%  "

# Row 42: HuobiToken
$ws.Range("E42").Value = "  -0.11%  "

# Row 43: RenderToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.17"
$ws.Range("E43").Value = "  +4.58%  "

# Row 44: FraxShare
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.02"
$ws.Range("E44").Value = "  +2.88%  "

# Row 45: Kaspa
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0499"
$ws.Range("E45").Value = "  -1.60%  "

# Row 46: RocketPoolETH
$ws.Range("D46").Value = "1.975.18"
$ws.Range("E46").Value = "  +1.46%  "

# Row 47: Quant
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.84"
$ws.Range("E47").Value = "  +8.01%  "

# Row 48: WEMIXToken
$ws.Range("E48").Value = "  -1.38%  "

# Row 49: InjectiveProtocol
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.996"
$ws.Range("E49").Value = "  -0.29%  "

# Row 50: PaxDollar
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.91"
$ws.Range("E50").Value = "  -4.79%  "

# Row 51: BabyDogeCoin
$ws.Range("D51").Value = "0.0₆0125"
$ws.Range("E51").Value = "  +6.29%  "
